$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 233.33333
$ws.Range("I2").Value = 233.33333
$ws.Range("K2").Value = 233.33333
$ws.Range("M2").Value = -120.33333
$ws.Range("H4").Value = 203.22223
$ws.Range("I4").Value = 147.14285
$ws.Range("K4").Value = 147.14285
$ws.Range("M4").Value = -33.14285000000001
$ws.Range("H38").Value = 1313.1875
$ws.Range("I38").Value = 80.583336
$ws.Range("J38").Value = 5011
$ws.Range("K38").Value = 241.750008
$ws.Range("L38").Value = 15033
$ws.Range("M38").Value = 130.249992
$ws.Range("N38").Value = -15777
$ws.Range("H40").Value = 4876.4443
$ws.Range("I40").Value = 4882.25
$ws.Range("J40").Value = 4830
$ws.Range("K40").Value = 4882.25
$ws.Range("L40").Value = 4830
$ws.Range("M40").Value = -4707.25
$ws.Range("N40").Value = -5180
$ws.Range("H41").Value = 269.8846
$ws.Range("I41").Value = 173.76923
$ws.Range("J41").Value = 366
$ws.Range("K41").Value = 173.76923
$ws.Range("L41").Value = 366
$ws.Range("M41").Value = 266.23077
$ws.Range("N41").Value = -1246
$ws.Range("H55").Value = 403.5
$ws.Range("I55").Value = 133.22223
$ws.Range("J55").Value = 890
$ws.Range("K55").Value = 133.22223
$ws.Range("L55").Value = 890
$ws.Range("M55").Value = 80.77777
$ws.Range("N55").Value = -1318
$ws.Range("H58").Value = 80.333336
$ws.Range("I58").Value = 80.333336
$ws.Range("K58").Value = 241.000008
$ws.Range("M58").Value = -91.00000800000001
$ws.Range("H70").Value = 42499.2
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 52499
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 157497
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -158037
$ws.Range("H73").Value = 42499.2
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 52499
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 157497
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -159369
$ws.Range("H80").Value = 1443.4445
$ws.Range("J80").Value = 1487.6
$ws.Range("L80").Value = 4462.799999999999
$ws.Range("N80").Value = -6458.799999999999
$ws.Range("H83").Value = 1443.4445
$ws.Range("J83").Value = 1487.6
$ws.Range("L83").Value = 13388.4
$ws.Range("N83").Value = -23372.4
$ws.Range("H92").Value = 1036.25
$ws.Range("I92").Value = 970.4286
$ws.Range("K92").Value = 970.4286
$ws.Range("M92").Value = 277.5714
$ws.Range("H101").Value = 1486.7142
$ws.Range("I101").Value = 553
$ws.Range("K101").Value = 1659
$ws.Range("M101").Value = -37
$ws.Range("H106").Value = 3500.8823
$ws.Range("J106").Value = 4947.5
$ws.Range("L106").Value = 4947.5
$ws.Range("N106").Value = -6209.5
$ws.Range("H116").Value = 5937.3
$ws.Range("I116").Value = 7914.3335
$ws.Range("J116").Value = 2971.75
$ws.Range("K116").Value = 7914.3335
$ws.Range("L116").Value = 2971.75
$ws.Range("M116").Value = -4472.3335
$ws.Range("N116").Value = -9855.75
$ws.Range("H125").Value = 2033.4615
$ws.Range("I125").Value = 742.75
$ws.Range("K125").Value = 6684.75
$ws.Range("M125").Value = -4224.75
$ws.Range("H131").Value = 4122.5
$ws.Range("I131").Value = 1154.7368
$ws.Range("K131").Value = 3464.2104
$ws.Range("M131").Value = 1575.7896
$ws.Range("H132").Value = 5792.939
$ws.Range("I132").Value = 3101.8838
$ws.Range("K132").Value = 9305.651400000001
$ws.Range("M132").Value = -6775.651400000001
$ws.Range("H135").Value = 1448.3784
$ws.Range("I135").Value = 956.3226
$ws.Range("J135").Value = 3990.6667
$ws.Range("K135").Value = 8606.903399999999
$ws.Range("L135").Value = 35916.0003
$ws.Range("M135").Value = -6071.903399999999
$ws.Range("N135").Value = -40986.0003
$ws.Range("H137").Value = 4879.795
$ws.Range("I137").Value = 5648.393
$ws.Range("J137").Value = 2923.3635
$ws.Range("K137").Value = 16945.179
$ws.Range("L137").Value = 8770.0905
$ws.Range("M137").Value = -14395.179
$ws.Range("N137").Value = -13870.0905
$ws.Range("H138").Value = 2200.125
$ws.Range("I138").Value = 1124.7587
$ws.Range("K138").Value = 3374.2761
$ws.Range("M138").Value = 1765.7239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5512.104
$ws.Range("I32").Value = 5012.427
$ws.Range("K32").Value = 5012.427
$ws.Range("M32").Value = -4725.427
$ws.Range("H51").Value = 27564
$ws.Range("J51").Value = 27564
$ws.Range("L51").Value = 27564
$ws.Range("N51").Value = -29076
$ws.Range("H61").Value = 5266.5
$ws.Range("I61").Value = 5110.5
$ws.Range("J61").Value = 6124.5
$ws.Range("K61").Value = 5110.5
$ws.Range("L61").Value = 6124.5
$ws.Range("M61").Value = -4898.5
$ws.Range("N61").Value = -6548.5
$ws.Range("H74").Value = 2916.3125
$ws.Range("I74").Value = 2060.2727
$ws.Range("J74").Value = 4799.6
$ws.Range("K74").Value = 2060.2727
$ws.Range("L74").Value = 4799.6
$ws.Range("M74").Value = -1186.2727
$ws.Range("N74").Value = -6547.6
$ws.Range("H77").Value = 2916.3125
$ws.Range("I77").Value = 2060.2727
$ws.Range("J77").Value = 4799.6
$ws.Range("K77").Value = 10301.3635
$ws.Range("L77").Value = 23998
$ws.Range("M77").Value = -5933.363499999999
$ws.Range("N77").Value = -32734
$ws.Range("H110").Value = 2919.75
$ws.Range("I110").Value = 3593.3333
$ws.Range("K110").Value = 3593.3333
$ws.Range("M110").Value = -1548.3333
$ws.Range("H122").Value = 36440.1
$ws.Range("I122").Value = 4468.3335
$ws.Range("J122").Value = 50142.285
$ws.Range("K122").Value = 13405.0005
$ws.Range("L122").Value = 150426.855
$ws.Range("M122").Value = -10955.0005
$ws.Range("N122").Value = -155326.855
$ws.Range("H132").Value = 2714.913
$ws.Range("I132").Value = 2582.3
$ws.Range("J132").Value = 3599
$ws.Range("K132").Value = 7746.900000000001
$ws.Range("L132").Value = 10797
$ws.Range("M132").Value = -5216.900000000001
$ws.Range("N132").Value = -15857
$ws.Range("H135").Value = 102898.336
$ws.Range("I135").Value = 101345
$ws.Range("J135").Value = 103675
$ws.Range("K135").Value = 101345
$ws.Range("L135").Value = 103675
$ws.Range("M135").Value = -96275
$ws.Range("N135").Value = -113815
$ws.Range("H136").Value = 5266.5
$ws.Range("I136").Value = 5110.5
$ws.Range("J136").Value = 6124.5
$ws.Range("K136").Value = 15331.5
$ws.Range("L136").Value = 18373.5
$ws.Range("M136").Value = -12781.5
$ws.Range("N136").Value = -23473.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H20").Value = 2531.1052
$ws.Range("I20").Value = 2280.2
$ws.Range("J20").Value = 3472
$ws.Range("K20").Value = 2280.2
$ws.Range("L20").Value = 3472
$ws.Range("M20").Value = -2033.2
$ws.Range("N20").Value = -3966
$ws.Range("H75").Value = 14800
$ws.Range("I75").Value = 14800
$ws.Range("K75").Value = 14800
$ws.Range("M75").Value = -13864
$ws.Range("H78").Value = 14800
$ws.Range("I78").Value = 14800
$ws.Range("K78").Value = 44400
$ws.Range("M78").Value = -39720
$ws.Range("H86").Value = 2616.5
$ws.Range("I86").Value = 3050
$ws.Range("J86").Value = 1749.5
$ws.Range("K86").Value = 3050
$ws.Range("L86").Value = 1749.5
$ws.Range("M86").Value = -1927
$ws.Range("N86").Value = -3995.5
$ws.Range("H89").Value = 2616.5
$ws.Range("I89").Value = 3050
$ws.Range("J89").Value = 1749.5
$ws.Range("K89").Value = 15250
$ws.Range("L89").Value = 8747.5
$ws.Range("M89").Value = -9634
$ws.Range("N89").Value = -19979.5
$ws.Range("H94").Value = 847.375
$ws.Range("J94").Value = 1254.75
$ws.Range("L94").Value = 1254.75
$ws.Range("N94").Value = -2156.75
$ws.Range("H105").Value = 2127.9285
$ws.Range("I105").Value = 2060.8462
$ws.Range("K105").Value = 2060.8462
$ws.Range("M105").Value = -313.8462
$ws.Range("H107").Value = 4951.381
$ws.Range("I107").Value = 4362.706
$ws.Range("K107").Value = 4362.706
$ws.Range("M107").Value = -2442.706
$ws.Range("H132").Value = 154123.75
$ws.Range("J132").Value = 154123.75
$ws.Range("L132").Value = 154123.75
$ws.Range("N132").Value = -164243.75
$ws.Range("H134").Value = 2896.087
$ws.Range("I134").Value = 2830.122
$ws.Range("J134").Value = 3437
$ws.Range("K134").Value = 8490.366
$ws.Range("L134").Value = 10311
$ws.Range("M134").Value = -5955.366
$ws.Range("N134").Value = -15381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1691666.6
$ws.Range("I6").Value = 1691666.6
$ws.Range("K6").Value = 1691666.6
$ws.Range("M6").Value = -1691553.6
$ws.Range("H16").Value = 2403.25
$ws.Range("I16").Value = 1601.2727
$ws.Range("J16").Value = 4167.6
$ws.Range("K16").Value = 1601.2727
$ws.Range("L16").Value = 4167.6
$ws.Range("M16").Value = -1314.2727
$ws.Range("N16").Value = -4741.6
$ws.Range("H31").Value = 1727.32
$ws.Range("I31").Value = 1438.7894
$ws.Range("J31").Value = 2641
$ws.Range("K31").Value = 1438.7894
$ws.Range("L31").Value = 2641
$ws.Range("M31").Value = -1143.7894
$ws.Range("N31").Value = -3231
$ws.Range("H34").Value = 1727.32
$ws.Range("I34").Value = 1438.7894
$ws.Range("J34").Value = 2641
$ws.Range("K34").Value = 1438.7894
$ws.Range("L34").Value = 2641
$ws.Range("M34").Value = -1236.7894
$ws.Range("N34").Value = -3045
$ws.Range("H56").Value = 49832.168
$ws.Range("I56").Value = 54664.332
$ws.Range("J56").Value = 45000
$ws.Range("K56").Value = 54664.332
$ws.Range("L56").Value = 45000
$ws.Range("M56").Value = -53819.332
$ws.Range("N56").Value = -46690
$ws.Range("H58").Value = 2413.4285
$ws.Range("I58").Value = 973.5
$ws.Range("J58").Value = 4333.3335
$ws.Range("K58").Value = 973.5
$ws.Range("L58").Value = 4333.3335
$ws.Range("M58").Value = -770.5
$ws.Range("N58").Value = -4739.3335
$ws.Range("H99").Value = 11474.448
$ws.Range("J99").Value = 17399.166
$ws.Range("L99").Value = 17399.166
$ws.Range("N99").Value = -20395.166
$ws.Range("H105").Value = 2284.5715
$ws.Range("I105").Value = 999
$ws.Range("J105").Value = 2498.8333
$ws.Range("K105").Value = 999
$ws.Range("L105").Value = 2498.8333
$ws.Range("M105").Value = 748
$ws.Range("N105").Value = -5992.8333
$ws.Range("H107").Value = 1159.5358
$ws.Range("I107").Value = 1051.3
$ws.Range("K107").Value = 1051.3
$ws.Range("M107").Value = 868.7
$ws.Range("H113").Value = 2403.25
$ws.Range("I113").Value = 1601.2727
$ws.Range("J113").Value = 4167.6
$ws.Range("K113").Value = 1601.2727
$ws.Range("L113").Value = 4167.6
$ws.Range("M113").Value = 568.7273
$ws.Range("N113").Value = -8507.6
$ws.Range("H122").Value = 9021.799999999999
$ws.Range("I122").Value = 9626.166999999999
$ws.Range("K122").Value = 28878.501
$ws.Range("M122").Value = -26428.501
$ws.Range("H126").Value = 11474.448
$ws.Range("J126").Value = 17399.166
$ws.Range("L126").Value = 52197.49800000001
$ws.Range("N126").Value = -57137.49800000001
$ws.Range("H132").Value = 2561.0967
$ws.Range("I132").Value = 2392.9656
$ws.Range("K132").Value = 7178.8968
$ws.Range("M132").Value = -4648.8968
$ws.Range("H134").Value = 2571.17
$ws.Range("I134").Value = 2161.7112
$ws.Range("K134").Value = 6485.133600000001
$ws.Range("M134").Value = -3950.133600000001
$ws.Range("H136").Value = 2413.4285
$ws.Range("I136").Value = 973.5
$ws.Range("J136").Value = 4333.3335
$ws.Range("K136").Value = 2920.5
$ws.Range("L136").Value = 13000.0005
$ws.Range("M136").Value = -370.5
$ws.Range("N136").Value = -18100.0005
$ws.Range("H138").Value = 78819.2
$ws.Range("J138").Value = 78819.2
$ws.Range("L138").Value = 78819.2
$ws.Range("N138").Value = -89099.2
$ws.Range("H141").Value = 123477.36
$ws.Range("J141").Value = 130269.914
$ws.Range("L141").Value = 130269.914
$ws.Range("N141").Value = -140629.914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 88.30303000000001
$ws.Range("J2").Value = 90.888885
$ws.Range("L2").Value = 545.33331
$ws.Range("N2").Value = -771.33331
$ws.Range("H17").Value = 330
$ws.Range("J17").Value = 400
$ws.Range("L17").Value = 1200
$ws.Range("N17").Value = -1538
$ws.Range("H34").Value = 379.72726
$ws.Range("J34").Value = 543.8333
$ws.Range("L34").Value = 1631.4999
$ws.Range("N34").Value = -1799.4999
$ws.Range("H39").Value = 9482.4
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 9482.4
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 28447.2
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -29035.2
$ws.Range("H55").Value = 3319.5293
$ws.Range("J55").Value = 8333.333000000001
$ws.Range("L55").Value = 24999.999
$ws.Range("N55").Value = -25353.999
$ws.Range("H56").Value = 11218.444
$ws.Range("I56").Value = 11218.444
$ws.Range("K56").Value = 11218.444
$ws.Range("M56").Value = -10688.444
$ws.Range("H86").Value = 473.76923
$ws.Range("J86").Value = 591.8
$ws.Range("L86").Value = 1775.4
$ws.Range("N86").Value = -4147.4
$ws.Range("H89").Value = 473.76923
$ws.Range("J89").Value = 591.8
$ws.Range("L89").Value = 5326.2
$ws.Range("N89").Value = -17182.2
$ws.Range("H97").Value = 1998.5
$ws.Range("I97").Value = 2746.75
$ws.Range("J97").Value = 1250.25
$ws.Range("K97").Value = 8240.25
$ws.Range("L97").Value = 3750.75
$ws.Range("M97").Value = -7744.25
$ws.Range("N97").Value = -4742.75
$ws.Range("H98").Value = 1216.4
$ws.Range("I98").Value = 1272.375
$ws.Range("K98").Value = 3817.125
$ws.Range("M98").Value = -2319.125
$ws.Range("H107").Value = 858.6
$ws.Range("I107").Value = 1075
$ws.Range("J107").Value = 804.5
$ws.Range("K107").Value = 3225
$ws.Range("L107").Value = 2413.5
$ws.Range("M107").Value = -1305
$ws.Range("N107").Value = -6253.5
$ws.Range("H109").Value = 1999.3125
$ws.Range("I109").Value = 1130.5
$ws.Range("K109").Value = 3391.5
$ws.Range("M109").Value = -2351.5
$ws.Range("H110").Value = 12255.5
$ws.Range("I110").Value = 12255.5
$ws.Range("K110").Value = 36766.5
$ws.Range("M110").Value = -32676.5
$ws.Range("H111").Value = 7000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 7000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 21000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -27134
$ws.Range("H112").Value = 3014.5
$ws.Range("I112").Value = 1999
$ws.Range("J112").Value = 4030
$ws.Range("K112").Value = 5997
$ws.Range("L112").Value = 12090
$ws.Range("M112").Value = -4889
$ws.Range("N112").Value = -14306
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H116").Value = 2792.1667
$ws.Range("I116").Value = 2563.5
$ws.Range("J116").Value = 3249.5
$ws.Range("K116").Value = 7690.5
$ws.Range("L116").Value = 9748.5
$ws.Range("M116").Value = -4248.5
$ws.Range("N116").Value = -16632.5
$ws.Range("H120").Value = 855.6667
$ws.Range("I120").Value = 855.6667
$ws.Range("K120").Value = 2567.0001
$ws.Range("M120").Value = 2270.9999
$ws.Range("H125").Value = 9933
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 9933
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29799
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -39639
$ws.Range("H131").Value = 542102.5
$ws.Range("I131").Value = 2972333
$ws.Range("J131").Value = 2051.2666
$ws.Range("K131").Value = 8916999
$ws.Range("L131").Value = 6153.7998
$ws.Range("M131").Value = -8911959
$ws.Range("N131").Value = -16233.7998
$ws.Range("H139").Value = 1976.7142
$ws.Range("I139").Value = 1577.0834
$ws.Range("K139").Value = 4731.2502
$ws.Range("M139").Value = 408.7497999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57.833332
$ws.Range("I2").Value = 57.833332
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 57.833332
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 55.166668
$ws.Range("N2").ClearContents()
$ws.Range("H3").Value = 5226.316
$ws.Range("I3").Value = 1389.5
$ws.Range("J3").Value = 9489.444
$ws.Range("K3").Value = 1389.5
$ws.Range("L3").Value = 9489.444
$ws.Range("M3").Value = -1273.5
$ws.Range("N3").Value = -9721.444
$ws.Range("H11").Value = 617558.8
$ws.Range("I11").Value = 653625
$ws.Range("J11").Value = 40500
$ws.Range("K11").Value = 653625
$ws.Range("L11").Value = 40500
$ws.Range("M11").Value = -653486
$ws.Range("N11").Value = -40778
$ws.Range("H14").Value = 422286.5
$ws.Range("I14").Value = 574551.5600000001
$ws.Range("J14").Value = 67001.336
$ws.Range("K14").Value = 574551.5600000001
$ws.Range("L14").Value = 67001.336
$ws.Range("M14").Value = -574383.5600000001
$ws.Range("N14").Value = -67337.336
$ws.Range("H80").Value = 2659.8262
$ws.Range("I80").Value = 2226.1
$ws.Range("J80").Value = 2993.4614
$ws.Range("K80").Value = 2226.1
$ws.Range("L80").Value = 2993.4614
$ws.Range("M80").Value = -1228.1
$ws.Range("N80").Value = -4989.4614
$ws.Range("H83").Value = 2659.8262
$ws.Range("I83").Value = 2226.1
$ws.Range("J83").Value = 2993.4614
$ws.Range("K83").Value = 11130.5
$ws.Range("L83").Value = 14967.307
$ws.Range("M83").Value = -6138.5
$ws.Range("N83").Value = -24951.307
$ws.Range("H113").Value = 2175
$ws.Range("I113").Value = 1901.8572
$ws.Range("K113").Value = 1901.8572
$ws.Range("M113").Value = 268.1428000000001
$ws.Range("H132").Value = 2304.4211
$ws.Range("I132").Value = 1461.1538
$ws.Range("J132").Value = 4131.5
$ws.Range("K132").Value = 4383.4614
$ws.Range("L132").Value = 12394.5
$ws.Range("M132").Value = -1853.4614
$ws.Range("N132").Value = -17454.5
$ws.Range("H140").Value = 145000
$ws.Range("J140").Value = 145000
$ws.Range("L140").Value = 145000
$ws.Range("N140").Value = -155360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2919.5833
$ws.Range("I7").Value = 2480.625
$ws.Range("J7").Value = 3797.5
$ws.Range("K7").Value = 2480.625
$ws.Range("L7").Value = 3797.5
$ws.Range("M7").Value = -2368.625
$ws.Range("N7").Value = -4021.5
$ws.Range("H16").Value = 712.9
$ws.Range("I16").Value = 687.2632
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 687.2632
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -517.2632
$ws.Range("N16").Value = -1540
$ws.Range("H22").Value = 995
$ws.Range("I22").Value = 995
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 995
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -700
$ws.Range("N22").Value = -1585
$ws.Range("H27").Value = 995
$ws.Range("I27").Value = 995
$ws.Range("J27").Value = 995
$ws.Range("K27").Value = 995
$ws.Range("L27").Value = 995
$ws.Range("M27").Value = -888
$ws.Range("N27").Value = -1209
$ws.Range("H40").Value = 7596
$ws.Range("I40").Value = 7596
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7596
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7460
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 2703.9565
$ws.Range("J46").Value = 3728
$ws.Range("L46").Value = 3728
$ws.Range("N46").Value = -4104
$ws.Range("H61").Value = 22063.6
$ws.Range("I61").Value = 22063.6
$ws.Range("K61").Value = 22063.6
$ws.Range("M61").Value = -21861.6
$ws.Range("H68").Value = 2972.25
$ws.Range("J68").Value = 3500
$ws.Range("L68").Value = 3500
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 2972.25
$ws.Range("J71").Value = 3500
$ws.Range("L71").Value = 17500
$ws.Range("N71").Value = -24988
$ws.Range("H74").Value = 46749.25
$ws.Range("I74").Value = 37000
$ws.Range("K74").Value = 37000
$ws.Range("M74").Value = -36002
$ws.Range("H77").Value = 46749.25
$ws.Range("I77").Value = 37000
$ws.Range("K77").Value = 111000
$ws.Range("M77").Value = -106008
$ws.Range("H82").Value = 2265
$ws.Range("J82").Value = 1600
$ws.Range("L82").Value = 1600
$ws.Range("N82").Value = -2322
$ws.Range("H85").Value = 2265
$ws.Range("J85").Value = 1600
$ws.Range("L85").Value = 1600
$ws.Range("N85").Value = -4096
$ws.Range("H100").Value = 2456.5
$ws.Range("I100").Value = 1726
$ws.Range("J100").Value = 2700
$ws.Range("K100").Value = 1726
$ws.Range("L100").Value = 2700
$ws.Range("M100").Value = -1185
$ws.Range("N100").Value = -3782
$ws.Range("H112").Value = 98875.5
$ws.Range("J112").Value = 98875.5
$ws.Range("L112").Value = 98875.5
$ws.Range("N112").Value = -101829.5
$ws.Range("H113").Value = 22063.6
$ws.Range("I113").Value = 22063.6
$ws.Range("K113").Value = 22063.6
$ws.Range("M113").Value = -19893.6
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H126").Value = 2919.5833
$ws.Range("I126").Value = 2480.625
$ws.Range("J126").Value = 3797.5
$ws.Range("K126").Value = 7441.875
$ws.Range("L126").Value = 11392.5
$ws.Range("M126").Value = -4971.875
$ws.Range("N126").Value = -16332.5
$ws.Range("H132").Value = 140004.5
$ws.Range("I132").Value = 170839.33
$ws.Range("K132").Value = 512517.99
$ws.Range("M132").Value = -509987.99
$ws.Range("H136").Value = 5353.9355
$ws.Range("I136").Value = 5841.269
$ws.Range("J136").Value = 2819.8
$ws.Range("K136").Value = 17523.807
$ws.Range("L136").Value = 8459.400000000001
$ws.Range("M136").Value = -14973.807
$ws.Range("N136").Value = -13559.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 48124.5
$ws.Range("J39").Value = 29750
$ws.Range("L39").Value = 29750
$ws.Range("N39").Value = -30576
$ws.Range("H45").Value = 14237.667
$ws.Range("J45").Value = 14237.667
$ws.Range("L45").Value = 14237.667
$ws.Range("N45").Value = -15219.667
$ws.Range("H54").Value = 41676.855
$ws.Range("I54").Value = 32500
$ws.Range("J54").Value = 45347.6
$ws.Range("K54").Value = 32500
$ws.Range("L54").Value = 45347.6
$ws.Range("M54").Value = -31980
$ws.Range("N54").Value = -46387.6
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 7549.5713
$ws.Range("I81").Value = 6739
$ws.Range("J81").Value = 7999.8887
$ws.Range("K81").Value = 13478
$ws.Range("L81").Value = 15999.7774
$ws.Range("M81").Value = -12417
$ws.Range("N81").Value = -18121.7774
$ws.Range("H84").Value = 7549.5713
$ws.Range("I84").Value = 6739
$ws.Range("J84").Value = 7999.8887
$ws.Range("K84").Value = 67390
$ws.Range("L84").Value = 79998.887
$ws.Range("M84").Value = -62086
$ws.Range("N84").Value = -90606.887
$ws.Range("H100").Value = 598.73334
$ws.Range("I100").Value = 598.73334
$ws.Range("K100").Value = 1197.46668
$ws.Range("M100").Value = -656.46668
$ws.Range("H107").Value = 1423.125
$ws.Range("J107").Value = 1657.5
$ws.Range("L107").Value = 4972.5
$ws.Range("N107").Value = -8812.5
$ws.Range("H113").Value = 4991
$ws.Range("I113").Value = 4991.3335
$ws.Range("K113").Value = 14974.0005
$ws.Range("M113").Value = -12804.0005
$ws.Range("H124").Value = 86661.5
$ws.Range("J124").Value = 86661.5
$ws.Range("L124").Value = 86661.5
$ws.Range("N124").Value = -96481.5
$ws.Range("H126").Value = 3470.56
$ws.Range("I126").Value = 4407.3125
$ws.Range("J126").Value = 1805.2222
$ws.Range("K126").Value = 13221.9375
$ws.Range("L126").Value = 5415.6666
$ws.Range("M126").Value = -10751.9375
$ws.Range("N126").Value = -10355.6666
$ws.Range("H132").Value = 3344.7144
$ws.Range("I132").Value = 3093.5945
$ws.Range("J132").Value = 5203
$ws.Range("K132").Value = 9280.783500000001
$ws.Range("L132").Value = 15609
$ws.Range("M132").Value = -6750.783500000001
$ws.Range("N132").Value = -20669
$ws.Range("H136").Value = 1970.9778
$ws.Range("I136").Value = 1791.8918
$ws.Range("K136").Value = 5375.6754
$ws.Range("M136").Value = -2825.6754
$ws.Range("H139").Value = 69900
$ws.Range("J139").Value = 69900
$ws.Range("L139").Value = 69900
$ws.Range("N139").Value = -80180
